# Fix printout layout: shrink the title line ("Gebäude #g / Raum #r")
# from 36pt (72 half-points) down to 32pt (64 half-points), and bump the
# paragraph mark's own size from 30pt (60 half-points) up to match (32pt).
#
# This affects the whole first paragraph (the drawn separator line run plus
# the four text runs "Gebäude ", "#g", " / Raum", " #r"), including the
# paragraph mark formatting stored in w:pPr/w:rPr.

$d = $word.ActiveDocument
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range

# Setting Font.Size / Font.SizeBi on the paragraph's Range updates both the
# ascii/east-asian run size (w:sz) and the complex-script size (w:szCs) for
# every run in the paragraph, as well as the paragraph mark's own rPr.
$titleRange.Font.Size = 32
$titleRange.Font.SizeBi = 32
